# Update crypto price/volume table per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.083.27'
$ws.Range('E2').Value = '  +5.55%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.261.50'
$ws.Range('E3').Value = '  +2.19%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.13'
$ws.Range('E5').Value = '  +3.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.97'
$ws.Range('E6').Value = '  +6.77%  '
$ws.Range('E7').Value = '  +3.67%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  +3.53%  '
$ws.Range('B10').Value = 'Avalanche'
$ws.Range('C10').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.75'
$ws.Range('E10').Value = '  +7.93%  '
$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.76'
$ws.Range('E11').Value = '  +9.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0800'
$ws.Range('E12').Value = '  +2.60%  '
$ws.Range('E13').Value = '  +3.38%  '
$ws.Range('E14').Value = '  +3.81%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.611.72'
$ws.Range('E15').Value = '  +2.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.16'
$ws.Range('E16').Value = '  +3.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.260.63'
$ws.Range('E17').Value = '  +2.34%  '
$ws.Range('E18').Value = '  +3.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.960.93'
$ws.Range('E19').Value = '  +5.41%  '
$ws.Range('E20').Value = '  +9.61%  '
$ws.Range('E21').Value = '  +2.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.96'
$ws.Range('E22').Value = '  +3.85%  '
$ws.Range('E23').Value = '  +2.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '241.95'
$ws.Range('E24').Value = '  +2.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.59'
$ws.Range('E25').Value = '  +5.64%  '
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.90'
$ws.Range('E27').Value = '  +4.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.01'
$ws.Range('E28').Value = '  +3.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.31'
$ws.Range('E29').Value = '  +12.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.69'
$ws.Range('E30').Value = '  +5.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.35'
$ws.Range('E31').Value = '  +7.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '159.40'
$ws.Range('E32').Value = '  +1.39%  '
$ws.Range('E34').Value = '  +4.22%  '
$ws.Range('E35').Value = '  +4.88%  '
$ws.Range('E36').Value = '  +4.45%  '
$ws.Range('E37').Value = '  +2.97%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '16.67'
$ws.Range('E38').Value = '  +9.40%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.105'
$ws.Range('E39').Value = '  +6.45%  '
$ws.Range('E40').Value = '  +4.37%  '
$ws.Range('E41').Value = '  +5.02%  '
$ws.Range('E42').Value = '  +5.90%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.054.20'
$ws.Range('E43').Value = '  -2.77%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.97'
$ws.Range('E44').Value = '  +12.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0280'
$ws.Range('E45').Value = '  +3.91%  '
$ws.Range('E46').Value = '  +0.94%  '
$ws.Range('E47').Value = '  +8.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.01'
$ws.Range('E48').Value = '  -3.71%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.483.41'
$ws.Range('E49').Value = '  +2.45%  '
$ws.Range('E50').Value = '  +3.49%  '
$ws.Range('E51').Value = '  +4.63%  '
